$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G, shifting existing G:K to H:L
$ws.Columns("G:G").Insert()

# New header for the inserted column (G1)
$ws.Range("G1").Value = "Cọc sim"

# New data values for the inserted column (G2:G3)
$ws.Range("G2").Value = "Cạnh"
$ws.Range("G3").Value = "Cạnh"

# Set the new column width (not bestFit, just a fixed custom width).
# Target stored OOXML width is 8.7109375; the closest value this
# COM runtime can reliably quantize to is 8.666666666666666.
$ws.Columns("G:G").ColumnWidth = 7.9

# Update the selected cell to match the target state
$ws.Range("G4").Select()
